$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Remove the "Phụ cấp tại SÓC TRĂNG" row first (row 26), then the
# "Phụ cấp tại LONG XUYÊN" row (row 14), so row indices stay valid
# for each delete operation.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(14).Delete()
